$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.463.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4837"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2897"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.904.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07420"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.178"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6624"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.417.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007776"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.138.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.378"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.98%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.204"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.388"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.945"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.446"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.332"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09204"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.046"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05077"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7571"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.709"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01877"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.652"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9175"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.089"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.982"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4355"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.670"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.593"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.14%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1329"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -12.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.904"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05709"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.34%  "
